$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be determined.`n"
$ws.Range("D2").Value = "no_decision, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been executed successfully.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: No decision was made regarding the movie to be shown on Friday.`n"
$ws.Range("D4").Value = "no_decision, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire rights for both movies. The committee has decided to show both `"Oppenheimer`" and `"Barbie`" on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been selected for acquisition.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie selection for Friday.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" have been acquired for Friday's screening.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision about the movie for Friday has not been made, so I am calling the no_decision function.`n"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: I have called the ``no_decision`` function, indicating that a clear decision about which movie to show on Friday could not be reached.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded, and there is no movie selected for Friday.`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision process ended without a plan for Friday's movie, resulting in no decision being made.`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The rights for both movies have been acquired.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has resulted in no agreement.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The committee did not arrive at a decision about which movie to show on Friday, so there will be no movie selected for acquisition.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The function has been executed, indicating that there was no decision made about the movie for Friday.`n"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision to acquire the rights for the movie `"Barbie`" has been recorded successfully.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision about which movie will be shown on Friday could not be made.`n"
$ws.Range("D31").Value = "no_decision, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" for Friday's screening.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached on the movie selection for Friday.`n"
$ws.Range("D33").Value = "no_decision, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The function has been successfully called, indicating that `"Barbie`" has been selected as the movie to be shown on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday remains unresolved, as there was no consensus reached in the discussion.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no specific movie was selected for Friday's showing.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for showing on Friday.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday's showing.`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to show on Friday.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision process concluded without an agreement on a movie for Friday.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision process did not result in an agreement on which movie to show on Friday.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision`" based on the discussion.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie can be made.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision process concluded without a clear choice for Friday's movie.`n"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday's showing.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D49").Value = "no_decision, "
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been made.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on what movie will be shown on Friday.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that the committee did not come to an agreement on a movie to show on Friday.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both `"Barbie`" and `"Oppenheimer`" for the movie to be shown on Friday.`n"
$ws.Range("D55").Value = "both_movies, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no consensus, so no movie has been chosen.`n"
$ws.Range("D56").Value = "no_decision, "
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for Friday's showing.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday cannot be made at this time.`n"
$ws.Range("D59").Value = "no_decision, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision about Friday's movie was not made, so there is no action to take.`n"
$ws.Range("D60").Value = "no_decision, "
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D62").Value = "both_movies, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie for Friday.`n"
$ws.Range("D63").Value = "no_decision, "
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made.`n"
$ws.Range("D64").Value = "no_decision, "
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached, so I will note that there was no decision made.`n"
$ws.Range("D65").Value = "no_decision, "
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("D66").Value = "no_decision, "
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been recorded, and there will be no movie shown on Friday as no agreement was reached.`n"
$ws.Range("D67").Value = "no_decision, "
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision about the movie to show on Friday resulted in no agreement. Therefore, I will call the no_decision function.`n"
$ws.Range("D68").Value = "no_decision, "
